# Updates the cryptos list (prices / 1h volume %) on Sheet1.
# For cells whose new text would otherwise be auto-recognized by Excel as a
# number (e.g. "307.25"), we briefly force a Text number format, assign the
# value, then restore the "Normal" style so the cell keeps behaving like the
# surrounding (unstyled) cells while still storing the value as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.239.12"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "
$ws.Cells.Item(3, 4).Value = "1.906.54"
$ws.Cells.Item(3, 5).Value = "  -0.01%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "307.25"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.19%  "
$ws.Cells.Item(6, 5).Value = "  -0.03%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5255"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.40%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3811"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +1.06%  "
$ws.Cells.Item(9, 5).Value = "  +0.20%  "
$ws.Cells.Item(10, 5).Value = "  +2.47%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.9031"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.49%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.08186"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -3.51%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "96.23"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.82%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "5.363"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.08%  "
$ws.Cells.Item(15, 4).Value = "1.456.57"
$ws.Cells.Item(15, 5).Value = "  -23.67%  "
$ws.Cells.Item(16, 5).Value = "  -0.06%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.000008668"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "14.76"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.33%  "
$ws.Cells.Item(19, 5).Value = "  +0.01%  "
$ws.Cells.Item(20, 4).Value = "27.269.29"
$ws.Cells.Item(20, 5).Value = "  +0.09%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "5.118"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.34%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "6.504"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.87%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "150.15"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.00%  "
$ws.Cells.Item(25, 5).Value = "  -0.88%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "18.25"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.35%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "1.742"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.02%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "116.68"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.24%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "4.847"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.19%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "4.845"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.64%  "
$ws.Cells.Item(31, 5).Value = "  -0.87%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.8332"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +4.18%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.05058"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.23%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "1.228"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.36%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "2.990"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.39%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "2.732"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.40%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "3.334"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -3.02%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.5815"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.03%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.02008"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.15%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.076"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.02%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "9.184"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.04%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "6.610"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.38%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "117.46"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.22%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.1522"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.13%  "
$ws.Cells.Item(45, 5).Value = "  +1.27%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "10.19"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.13%  "
$ws.Cells.Item(47, 5).Value = "  -0.07%  "
$ws.Cells.Item(48, 5).Value = "  +0.88%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "38.87"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.98%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.06125"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.79%  "
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "64.48"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.54%  "
